$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text (matches original inlineStr formatting)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '26.158.66'
$ws.Cells.Item(2, 5).Value = '  -0.06%  '

$ws.Cells.Item(3, 4).Value = '1.656.07'
$ws.Cells.Item(3, 5).Value = '  -0.15%  '

$ws.Cells.Item(4, 5).Value = '  -0.18%  '

$ws.Cells.Item(5, 4).Value = '218.74'
$ws.Cells.Item(5, 5).Value = '  -0.22%  '

$ws.Cells.Item(6, 4).Value = '0.5237'
$ws.Cells.Item(6, 5).Value = '  +0.15%  '

$ws.Cells.Item(7, 5).Value = '  -0.18%  '

$ws.Cells.Item(9, 4).Value = '0.06348'
$ws.Cells.Item(9, 5).Value = '  +0.74%  '

$ws.Cells.Item(10, 4).Value = '20.53'
$ws.Cells.Item(10, 5).Value = '  -0.61%  '

$ws.Cells.Item(11, 4).Value = '0.07681'
$ws.Cells.Item(11, 5).Value = '  -1.76%  '

$ws.Cells.Item(12, 4).Value = '4.619'

$ws.Cells.Item(13, 4).Value = '1.652.66'
$ws.Cells.Item(13, 5).Value = '  -0.28%  '

$ws.Cells.Item(14, 4).Value = '1.884.24'
$ws.Cells.Item(14, 5).Value = '  -0.12%  '

$ws.Cells.Item(15, 4).Value = '0.5613'
$ws.Cells.Item(15, 5).Value = '  +1.11%  '

$ws.Cells.Item(16, 4).Value = '0.0₅8183'
$ws.Cells.Item(16, 5).Value = '  +2.01%  '

$ws.Cells.Item(17, 5).Value = '  +0.59%  '

$ws.Cells.Item(18, 4).Value = '26.148.63'
$ws.Cells.Item(18, 5).Value = '  -0.12%  '

$ws.Cells.Item(19, 4).Value = '1.003'
$ws.Cells.Item(19, 5).Value = '  -0.22%  '

$ws.Cells.Item(20, 4).Value = '4.652'
$ws.Cells.Item(20, 5).Value = '  +0.19%  '

$ws.Cells.Item(21, 5).Value = '  +3.51%  '

$ws.Cells.Item(22, 4).Value = '193.12'
$ws.Cells.Item(22, 5).Value = '  -1.53%  '

$ws.Cells.Item(23, 4).Value = '5.954'
$ws.Cells.Item(23, 5).Value = '  -0.21%  '

$ws.Cells.Item(24, 5).Value = '  -0.19%  '

$ws.Cells.Item(25, 4).Value = '145.22'
$ws.Cells.Item(25, 5).Value = '  -0.69%  '

$ws.Cells.Item(26, 4).Value = '0.1195'
$ws.Cells.Item(26, 5).Value = '  -0.57%  '

$ws.Cells.Item(27, 5).Value = '  +1.51%  '

$ws.Cells.Item(28, 5).Value = '  -0.35%  '

$ws.Cells.Item(29, 4).Value = '1.516'
$ws.Cells.Item(29, 5).Value = '  +1.61%  '

$ws.Cells.Item(30, 4).Value = '0.05480'
$ws.Cells.Item(30, 5).Value = '  -4.58%  '

$ws.Cells.Item(31, 4).Value = '1.273'
$ws.Cells.Item(31, 5).Value = '  -0.14%  '

$ws.Cells.Item(32, 4).Value = '3.465'
$ws.Cells.Item(32, 5).Value = '  -0.91%  '

$ws.Cells.Item(33, 4).Value = '3.366'
$ws.Cells.Item(33, 5).Value = '  -0.45%  '

$ws.Cells.Item(34, 4).Value = '1.559'
$ws.Cells.Item(34, 5).Value = '  -1.74%  '

$ws.Cells.Item(35, 4).Value = '0.9507'
$ws.Cells.Item(35, 5).Value = '  -0.63%  '

$ws.Cells.Item(36, 4).Value = '2.778'
$ws.Cells.Item(36, 5).Value = '  -0.93%  '

$ws.Cells.Item(37, 5).Value = '  -0.76%  '

$ws.Cells.Item(38, 4).Value = '0.5682'
$ws.Cells.Item(38, 5).Value = '  -0.65%  '

$ws.Cells.Item(39, 4).Value = '0.01588'
$ws.Cells.Item(39, 5).Value = '  -0.56%  '

$ws.Cells.Item(40, 4).Value = '5.865'
$ws.Cells.Item(40, 5).Value = '  -1.69%  '

$ws.Cells.Item(41, 5).Value = '  -0.20%  '

$ws.Cells.Item(42, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(42, 4).Value = '0.8322'
$ws.Cells.Item(42, 5).Value = '  -2.15%  '

$ws.Cells.Item(43, 2).Value = 'Maker'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(43, 4).Value = '1.027.56'
$ws.Cells.Item(43, 5).Value = '  -3.47%  '

$ws.Cells.Item(44, 4).Value = '100.93'
$ws.Cells.Item(44, 5).Value = '  -2.92%  '

$ws.Cells.Item(45, 4).Value = '1.795.10'
$ws.Cells.Item(45, 5).Value = '  -0.11%  '

$ws.Cells.Item(46, 4).Value = '57.94'
$ws.Cells.Item(46, 5).Value = '  -0.09%  '

$ws.Cells.Item(47, 2).Value = 'Frax'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(47, 4).Value = '1.000'
$ws.Cells.Item(47, 5).Value = '  -0.54%  '

$ws.Cells.Item(48, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(48, 4).Value = '0.0₈104'
$ws.Cells.Item(48, 5).Value = '  +0.48%  '

$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(49, 4).Value = '8.056'
$ws.Cells.Item(49, 5).Value = '  +0.07%  '

$ws.Cells.Item(50, 2).Value = 'Mantle'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(50, 4).Value = '0.4347'
$ws.Cells.Item(50, 5).Value = '  -1.27%  '

$ws.Cells.Item(51, 4).Value = '0.05209'
$ws.Cells.Item(51, 5).Value = '  +0.18%  '
